$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("VT-AuthCapCredit-Generic")
$ws.Range("B2").Value = "Thu Aug 28 07:57:30 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 07:58:42 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 07:59:48 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 08:01:26 IST 2025"
$ws.Range("B6").Value = "Thu Aug 28 08:02:40 IST 2025"
$ws.Range("B7").Value = "Thu Aug 28 08:03:55 IST 2025"

$ws = $wb.Worksheets.Item("VT-AuthCapVoid-Generic")
$ws.Range("B2").Value = "Thu Aug 28 08:05:17 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 08:06:17 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 08:07:31 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 08:08:40 IST 2025"
$ws.Range("B6").Value = "Thu Aug 28 08:09:44 IST 2025"
$ws.Range("B7").Value = "Thu Aug 28 08:10:52 IST 2025"

$ws = $wb.Worksheets.Item("VT-ManualAuthCapture-Generic")
$ws.Range("B2").Value = "Thu Aug 28 08:41:33 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 08:42:24 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 08:43:11 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 08:44:00 IST 2025"
$ws.Range("B6").Value = "Thu Aug 28 08:44:49 IST 2025"
$ws.Range("B7").Value = "Thu Aug 28 08:45:40 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleCredit-DualCF-Generic")
$ws.Range("B2").Value = "Thu Aug 28 08:46:38 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 08:47:29 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 08:48:18 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 08:49:06 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleCredit-SingleCF-Generic")
$ws.Range("B2").Value = "Thu Aug 28 08:49:55 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 08:50:45 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 08:51:39 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 08:52:26 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleVoid-DualCF-Generic")
$ws.Range("B2").Value = "Thu Aug 28 08:53:11 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 08:54:01 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 08:54:44 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 08:55:35 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleVoid-NoCF-Generic")
$ws.Range("B2").Value = "Thu Aug 28 08:56:25 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 08:57:19 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 08:58:07 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 08:59:01 IST 2025"

$ws = $wb.Worksheets.Item("VT-SaleVoid-SingleCF-Generic")
$ws.Range("B2").Value = "Thu Aug 28 08:59:49 IST 2025"
$ws.Range("B3").Value = "Thu Aug 28 09:00:42 IST 2025"
$ws.Range("B4").Value = "Thu Aug 28 09:01:32 IST 2025"
$ws.Range("B5").Value = "Thu Aug 28 09:02:15 IST 2025"
